$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1662407049093986
$ws.Range("D2").Value = 0.2406812144209241
$ws.Range("E2").Value = 0.1632400434247998
$ws.Range("F2").Value = 0.9330991215533544
$ws.Range("G2").Value = 0.4657766087619706
$ws.Range("H2").Value = 0.5615101234332656
$ws.Range("I2").Value = 0.3661632451992496
$ws.Range("J2").Value = 0.151185536860881
$ws.Range("M2").Value = 0.576508063319352
$ws.Range("N2").Value = 1.324766333497138
$ws.Range("O2").Value = 2.010535735687853

# Row 3
$ws.Range("B3").Value = 0.1471823725918
$ws.Range("D3").Value = 0.2420724027139087
$ws.Range("E3").Value = 0.1650069050550975
$ws.Range("F3").Value = 0.919647765470728
$ws.Range("G3").Value = 0.4488392227571012
$ws.Range("H3").Value = 0.5574299128864624
$ws.Range("I3").Value = 0.3658408523029735
$ws.Range("J3").Value = 0.1534384553728287
$ws.Range("M3").Value = 0.5180052233953774
$ws.Range("N3").Value = 1.282048421507199
$ws.Range("O3").Value = 1.965764170496982

# Row 4
$ws.Range("B4").Value = 0.1354552622694314
$ws.Range("D4").Value = 0.2430230415639514
$ws.Range("E4").Value = 0.1661751917371994
$ws.Range("F4").Value = 0.9119589159233499
$ws.Range("G4").Value = 0.4387577022019684
$ws.Range("H4").Value = 0.5552159724355761
$ws.Range("I4").Value = 0.3658603238777545
$ws.Range("J4").Value = 0.1549139033592084
$ws.Range("M4").Value = 0.4820230782338655
$ws.Range("N4").Value = 1.256212769797202
$ws.Range("O4").Value = 1.939546256634088

# Row 5
$ws.Range("B5").Value = 0.1306704766576701
$ws.Range("D5").Value = 0.2434347085586595
$ws.Range("E5").Value = 0.1666722470441808
$ws.Range("F5").Value = 0.9089690380368936
$ws.Range("G5").Value = 0.4347292803214486
$ws.Range("H5").Value = 0.5543870791450018
$ws.Range("I5").Value = 0.3659231569378534
$ws.Range("J5").Value = 0.1555383004255138
$ws.Range("M5").Value = 0.4673459374938034
$ws.Range("N5").Value = 1.245785044692767
$ws.Range("O5").Value = 1.929181961860706

# Row 6
$ws.Range("B6").Value = 0.1298756239416292
$ws.Range("D6").Value = 0.2435045322964484
$ws.Range("E6").Value = 0.1667560489044586
$ws.Range("F6").Value = 0.9084812308350507
$ws.Range("G6").Value = 0.434065185270228
$ws.Range("H6").Value = 0.5542538709802187
$ws.Range("I6").Value = 0.3659369126941847
$ws.Range("J6").Value = 0.1556433777741155
$ws.Range("M6").Value = 0.464907993360562
$ws.Range("N6").Value = 1.244059647114824
$ws.Range("O6").Value = 1.927480284820376

# Row 7
$ws.Range("B7").Value = 0.1353907562417334
$ws.Range("D7").Value = 0.2430284951339772
$ws.Range("E7").Value = 0.1661818103084318
$ws.Range("F7").Value = 0.9119180128364377
$ws.Range("G7").Value = 0.4387030502132916
$ws.Range("H7").Value = 0.5552044968159322
$ws.Range("I7").Value = 0.3658609486842792
$ws.Range("J7").Value = 0.1549222305443738
$ws.Range("M7").Value = 0.4818251927051591
$ws.Range("N7").Value = 1.256071728735549
$ws.Range("O7").Value = 1.939405185794101

# Row 8
$ws.Range("B8").Value = 0.1596748883226269
$ws.Range("D8").Value = 0.2411408926092591
$ws.Range("E8").Value = 0.163831936668279
$ws.Range("F8").Value = 0.9283426983590743
$ws.Range("G8").Value = 0.4598705033113788
$ws.Range("H8").Value = 0.5600428290006079
$ws.Range("I8").Value = 0.3660070918233771
$ws.Range("J8").Value = 0.1519431996187954
$ws.Range("M8").Value = 0.5563497138096238
$ws.Range("N8").Value = 1.309956510776715
$ws.Range("O8").Value = 1.994834391655331

# Row 9
$ws.Range("B9").Value = 0.2070784958907836
$ws.Range("D9").Value = 0.2382036222685855
$ws.Range("E9").Value = 0.159886190976156
$ws.Range("F9").Value = 0.9650810202261226
$ws.Range("G9").Value = 0.5039115515563424
$ws.Range("H9").Value = 0.5718413822429653
$ws.Range("I9").Value = 0.3680094933028819
$ws.Range("J9").Value = 0.1468337980814542
$ws.Range("M9").Value = 0.7019585033537794
$ws.Range("N9").Value = 1.418680917183224
$ws.Range("O9").Value = 2.113638542126807

# Row 10
$ws.Range("B10").Value = 0.2417532678818759
$ws.Range("D10").Value = 0.2365104022502322
$ws.Range("E10").Value = 0.1573915586016259
$ws.Range("F10").Value = 0.994844212469232
$ws.Range("G10").Value = 0.5378266232374074
$ws.Range("H10").Value = 0.5819189429804794
$ws.Range("I10").Value = 0.370514974748744
$ws.Range("J10").Value = 0.1435281773168864
$ws.Range("M10").Value = 0.8085573462053475
$ws.Range("N10").Value = 1.500347060190109
$ws.Range("O10").Value = 2.207118731935623

# Row 11
$ws.Range("B11").Value = 0.2574905359920194
$ws.Range("D11").Value = 0.2358408025943035
$ws.Range("E11").Value = 0.1563445906803658
$ws.Range("F11").Value = 1.00898853485738
$ws.Range("G11").Value = 0.5535971819498968
$ws.Range("H11").Value = 0.5868097174238045
$ws.Range("I11").Value = 0.3718769730888738
$ws.Range("J11").Value = 0.1421220671435179
$ws.Range("M11").Value = 0.8569584688636809
$ws.Range("N11").Value = 1.537871387701756
$ws.Range("O11").Value = 2.25099823365872

# Row 12
$ws.Range("B12").Value = 0.263444177063235
$ws.Range("D12").Value = 0.2356016968117203
$ws.Range("E12").Value = 0.1559607756535861
$ws.Range("F12").Value = 1.014431710911339
$ws.Range("G12").Value = 0.5596185055230336
$ws.Range("H12").Value = 0.5887057757092862
$ws.Range("I12").Value = 0.3724244726980857
$ws.Range("J12").Value = 0.1416036831868368
$ws.Range("M12").Value = 0.8752724520983577
$ws.Range("N12").Value = 1.552133174701652
$ws.Range("O12").Value = 2.267809454340977

# Row 13
$ws.Range("B13").Value = 0.2621622159605295
$ws.Range("D13").Value = 0.2356525498075257
$ws.Range("E13").Value = 0.1560428742667366
$ws.Range("F13").Value = 1.013255554255466
$ws.Range("G13").Value = 0.5583195088454715
$ws.Range("H13").Value = 0.5882954680688783
$ws.Range("I13").Value = 0.3723051506506749
$ws.Range("J13").Value = 0.1417146996241456
$ws.Range("M13").Value = 0.871328873260552
$ws.Range("N13").Value = 1.549059350816407
$ws.Range("O13").Value = 2.264180181508607

# Row 14
$ws.Range("B14").Value = 0.257980463010739
$ws.Range("D14").Value = 0.2358208415724476
$ws.Range("E14").Value = 0.1563127603913852
$ws.Range("F14").Value = 1.009434603801409
$ws.Range("G14").Value = 0.5540915700557179
$ws.Range("H14").Value = 0.5869648251065769
$ws.Range("I14").Value = 0.3719213813284625
$ws.Range("J14").Value = 0.1420791369819945
$ws.Range("M14").Value = 0.858465469055119
$ws.Range("N14").Value = 1.539043677922791
$ws.Range("O14").Value = 2.252377393713289

# Row 15
$ws.Range("B15").Value = 0.2554182594109875
$ws.Range("D15").Value = 0.2359258073746062
$ws.Range("E15").Value = 0.1564797211801405
$ws.Range("F15").Value = 1.007105496262653
$ws.Range("G15").Value = 0.5515082663783346
$ws.Range("H15").Value = 0.5861555005155736
$ws.Range("I15").Value = 0.3716904386364206
$ws.Range("J15").Value = 0.1423042001304324
$ws.Range("M15").Value = 0.8505843378411555
$ws.Range("N15").Value = 1.532915530568516
$ws.Range("O15").Value = 2.2451732459906

# Row 16
$ws.Range("B16").Value = 0.2407240258384888
$ws.Range("D16").Value = 0.2365561861418399
$ws.Range("E16").Value = 0.1574617495982142
$ws.Range("F16").Value = 0.9939320236524338
$ws.Range("G16").Value = 0.5368028804405327
$ws.Range("H16").Value = 0.5816054825242531
$ws.Range("I16").Value = 0.3704304183870448
$ws.Range("J16").Value = 0.1436220379652742
$ws.Range("M16").Value = 0.8053922675442635
$ws.Range("N16").Value = 1.497902142238559
$ws.Range("O16").Value = 2.204278382137147

# Row 17
$ws.Range("B17").Value = 0.2316999106313347
$ws.Range("D17").Value = 0.2369686710959726
$ws.Range("E17").Value = 0.1580867047413168
$ws.Range("F17").Value = 0.9860054856266913
$ws.Range("G17").Value = 0.5278693875654028
$ws.Range("H17").Value = 0.5788926539822512
$ws.Range("I17").Value = 0.3697142029288543
$ws.Range("J17").Value = 0.1444555234051137
$ws.Range("M17").Value = 0.7776440669996418
$ws.Range("N17").Value = 1.476517179070356
$ws.Range("O17").Value = 2.179537853048828

# Row 18
$ws.Range("B18").Value = 0.2265060749450072
$ws.Range("D18").Value = 0.2372153966329549
$ws.Range("E18").Value = 0.1584544305887832
$ws.Range("F18").Value = 0.9815032837498734
$ws.Range("G18").Value = 0.5227633023589391
$ws.Range("H18").Value = 0.5773611527269509
$ws.Range("I18").Value = 0.3693231861592992
$ws.Range("J18").Value = 0.1449441095066533
$ws.Range("M18").Value = 0.7616755300297484
$ws.Range("N18").Value = 1.464252412833076
$ws.Range("O18").Value = 2.165435273939522

# Row 19
$ws.Range("B19").Value = 0.2247469614999034
$ws.Range("D19").Value = 0.2373005615344042
$ws.Range("E19").Value = 0.158580355826528
$ws.Range("F19").Value = 0.9799886925012231
$ws.Range("G19").Value = 0.5210400005454972
$ws.Range("H19").Value = 0.5768475681950633
$ws.Range("I19").Value = 0.3691943962505704
$ws.Range("J19").Value = 0.1451111130813043
$ws.Range("M19").Value = 0.7562674381876207
$ws.Range("N19").Value = 1.460105886543914
$ws.Range("O19").Value = 2.160682279134022

# Row 20
$ws.Range("B20").Value = 0.2326608987972065
$ws.Range("D20").Value = 0.2369237808547524
$ws.Range("E20").Value = 0.1580193213592249
$ws.Range("F20").Value = 0.9868433864356234
$ws.Range("G20").Value = 0.5288170372694765
$ws.Range("H20").Value = 0.5791784540341922
$ws.Range("I20").Value = 0.3697882805466008
$ws.Range("J20").Value = 0.1443658463150239
$ws.Range("M20").Value = 0.7805988004632809
$ws.Range("N20").Value = 1.478790003518327
$ws.Range("O20").Value = 2.18215832561205

# Row 21
$ws.Range("B21").Value = 0.2592089042809391
$ws.Range("D21").Value = 0.2357710179759209
$ws.Range("E21").Value = 0.1562331448777243
$ws.Range("F21").Value = 1.010554547300785
$ws.Range("G21").Value = 0.5553320777652146
$ws.Range("H21").Value = 0.5873544725889985
$ws.Range("I21").Value = 0.3720332439840917
$ws.Range("J21").Value = 0.141971710465512
$ws.Range("M21").Value = 0.86224416620216
$ws.Range("N21").Value = 1.541984121825806
$ws.Range("O21").Value = 2.255838864824511

# Row 22
$ws.Range("B22").Value = 0.2765259947411209
$ws.Range("D22").Value = 0.2351018798655673
$ws.Range("E22").Value = 0.1551395107209945
$ws.Range("F22").Value = 1.026558448265305
$ws.Range("G22").Value = 0.5729489222825066
$ws.Range("H22").Value = 0.5929545856367326
$ws.Range("I22").Value = 0.3736853646768594
$ws.Range("J22").Value = 0.1404890886197627
$ws.Range("M22").Value = 0.9155191562828549
$ws.Range("N22").Value = 1.583588281711201
$ws.Range("O22").Value = 2.305130077566218

# Row 23
$ws.Range("B23").Value = 0.2672867663774525
$ws.Range("D23").Value = 0.2354513073415845
$ws.Range("E23").Value = 0.1557164518137455
$ws.Range("F23").Value = 1.017970433851971
$ws.Range("G23").Value = 0.563520115489581
$ws.Range("H23").Value = 0.5899422323240913
$ws.Range("I23").Value = 0.3727867493082471
$ws.Range("J23").Value = 0.1412728678924822
$ws.Range("M23").Value = 0.887093512950301
$ws.Range("N23").Value = 1.561356159312879
$ws.Range("O23").Value = 2.278718362981692

# Row 24
$ws.Range("B24").Value = 0.2322264539268986
$ws.Range("D24").Value = 0.2369440458836394
$ws.Range("E24").Value = 0.1580497591262784
$ws.Range("F24").Value = 0.9864644005635057
$ws.Range("G24").Value = 0.5283885117135725
$ws.Range("H24").Value = 0.5790491561699014
$ws.Range("I24").Value = 0.3697547254213589
$ws.Range("J24").Value = 0.1444063600368697
$ws.Range("M24").Value = 0.7792630143476487
$ws.Range("N24").Value = 1.477762366842484
$ws.Range("O24").Value = 2.180973232814523

# Row 25
$ws.Range("B25").Value = 0.1942799304984817
$ws.Range("D25").Value = 0.2389165262486372
$ws.Range("E25").Value = 0.1608826569196502
$ws.Range("F25").Value = 0.9546563525556877
$ws.Range("G25").Value = 0.4917247177410076
$ws.Range("H25").Value = 0.5684021383309528
$ws.Range("I25").Value = 0.3672853152251996
$ws.Range("J25").Value = 0.1481374327590608
$ws.Range("M25").Value = 0.6626302535026838
$ws.Range("N25").Value = 1.388948871434394
$ws.Range("O25").Value = 2.080413371653037
